$wb = $excel.ActiveWorkbook
$wsTCD = $wb.Worksheets.Item("TCD")
foreach ($pt in $wsTCD.PivotTables()) {
    Write-Host "pivot table: $($pt.Name)"
    try {
        $pt.RefreshTable()
        Write-Host "refresh ok"
    } catch {
        Write-Host "refresh err: $_"
    }
}
